$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Helper: swap the F:V contents (everything except the shared
# A/B/C/D/E "match metadata" columns) between two rows. The two
# rows describe the same match day / same two fixtures that had
# been paired with the wrong opponent, so the fix is a straight
# swap of columns F through V.
# ---------------------------------------------------------------
function Swap-RowData($row1, $row2) {
    for ($col = 6; $col -le 22; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

# Rows 122 / 123 : Mladost-Zeleznicar Pancevo vs Partizan-Vojvodina were
# swapped into the wrong slots.
Swap-RowData 122 123

# Rows 138 / 139 : Radnicki Nis-Vojvodina vs Mladost-Napredak likewise.
Swap-RowData 138 139

# ---------------------------------------------------------------
# Append the new match row 149 (Cukaricki 3-0 Mladost), copying the
# formatting of the last existing data row (148) so the new row's
# styling (bold/border on column A, datetime format on column E)
# matches the rest of the table.
# ---------------------------------------------------------------
$ws.Range("A148:V148").Copy()
$ws.Range("A149:V149").PasteSpecial(-4122)

$ws.Cells.Item(149, 1).Value2 = 148
$ws.Cells.Item(149, 2).Value2 = "serbia"
$ws.Cells.Item(149, 3).Value2 = "super-liga"
$ws.Cells.Item(149, 4).Value2 = "2023-2024"
$ws.Cells.Item(149, 5).Value2 = 45282.58333333334
$ws.Cells.Item(149, 6).Value2 = "Cukaricki"
$ws.Cells.Item(149, 7).Value2 = 3
$ws.Cells.Item(149, 8).Value2 = "Mladost"
$ws.Cells.Item(149, 9).Value2 = 0
$ws.Cells.Item(149, 10).Value2 = 1.36
$ws.Cells.Item(149, 11).Value2 = "27/09/2023 04:42"
$ws.Cells.Item(149, 12).Value2 = 1.42
$ws.Cells.Item(149, 13).Value2 = "22/12/2023 13:29"
$ws.Cells.Item(149, 14).Value2 = 4.41
$ws.Cells.Item(149, 15).Value2 = "27/09/2023 04:42"
$ws.Cells.Item(149, 16).Value2 = 4.67
$ws.Cells.Item(149, 17).Value2 = "22/12/2023 13:29"
$ws.Cells.Item(149, 18).Value2 = 6.61
$ws.Cells.Item(149, 19).Value2 = "27/09/2023 04:42"
$ws.Cells.Item(149, 20).Value2 = 6.58
$ws.Cells.Item(149, 21).Value2 = "22/12/2023 13:29"
$ws.Cells.Item(149, 22).Value2 = "https://www.betexplorer.com/football/serbia/super-liga/cukaricki-mladost-lucani/zXPec4Vb/"

Write-Output "done"
